# Natmi following Dr Hou advice
# Expand the Sending/Target cluster grid from 4x3 (ECs/FAPs/M2/sCs x ECs/FAPs/sCs)
# to a full 4x4 grid (adding M2 as a Target cluster) and refresh all the computed
# ligand/receptor expression metrics accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numCols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @{ Row = 2; A = "ECs"; B = "Gnas"; C = "Tshr"; D = "ECs"; Vals = @(3, 1, 111.005114, 333.015342, 0.07617677061978854, 0.07617677061978853, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 80.03302211729002, 720.2971990556101, 0.011680378759556, 0.01168037875955599) },
    @{ Row = 3; A = "ECs"; B = "Gnas"; C = "Tshr"; D = "FAPs"; Vals = @(3, 1, 111.005114, 333.015342, 0.07617677061978854, 0.07617677061978853, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 376.0080296709513, 3384.072267038562, 0.05487630089433133, 0.05487630089433131) },
    @{ Row = 4; A = "ECs"; B = "Gnas"; C = "Tshr"; D = "M2"; Vals = @(3, 1, 111.005114, 333.015342, 0.07617677061978854, 0.07617677061978853, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 12.95981005779534, 116.638290520158, 0.001891412895855775, 0.001891412895855775) },
    @{ Row = 5; A = "ECs"; B = "Gnas"; C = "Tshr"; D = "sCs"; Vals = @(3, 1, 111.005114, 333.015342, 0.07617677061978854, 0.07617677061978853, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 52.95628469336334, 476.6065622402701, 0.007728678070045447, 0.007728678070045444) },
    @{ Row = 6; A = "FAPs"; B = "Gnas"; C = "Tshr"; D = "ECs"; Vals = @(3, 1, 930.1503093333332, 2790.450928, 0.6383115534299679, 0.6383115534299678, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 670.6244207746934, 6035.61978697224, 0.09787394044144224, 0.09787394044144221) },
    @{ Row = 7; A = "FAPs"; B = "Gnas"; C = "Tshr"; D = "FAPs"; Vals = @(3, 1, 930.1503093333332, 2790.450928, 0.6383115534299679, 0.6383115534299678, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 3150.701553355934, 28356.31398020341, 0.4598275377829111, 0.459827537782911) },
    @{ Row = 8; A = "FAPs"; B = "Gnas"; C = "Tshr"; D = "M2"; Vals = @(3, 1, 930.1503093333332, 2790.450928, 0.6383115534299679, 0.6383115534299678, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 108.5947385645636, 977.352647081072, 0.01584880395832308, 0.01584880395832308) },
    @{ Row = 9; A = "FAPs"; B = "Gnas"; C = "Tshr"; D = "sCs"; Vals = @(3, 1, 930.1503093333332, 2790.450928, 0.6383115534299679, 0.6383115534299678, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 443.7390568210756, 3993.65151138968, 0.06476127124729156, 0.06476127124729153) },
    @{ Row = 10; A = "M2"; B = "Gnas"; C = "Tshr"; D = "ECs"; Vals = @(3, 1, 72.47344466666667, 217.420334, 0.04973458223193756, 0.04973458223193756, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 52.25226650299668, 470.27039852697, 0.00762593049286351, 0.007625930492863508) },
    @{ Row = 11; A = "M2"; B = "Gnas"; C = "Tshr"; D = "FAPs"; Vals = @(3, 1, 72.47344466666667, 217.420334, 0.04973458223193756, 0.04973458223193756, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 245.4895648553638, 2209.406083698274, 0.03582784984461771, 0.03582784984461771) },
    @{ Row = 12; A = "M2"; B = "Gnas"; C = "Tshr"; D = "M2"; Vals = @(3, 1, 72.47344466666667, 217.420334, 0.04973458223193756, 0.04973458223193756, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 8.461250507018445, 76.15125456316601, 0.001234872907293472, 0.001234872907293472) },
    @{ Row = 13; A = "M2"; B = "Gnas"; C = "Tshr"; D = "sCs"; Vals = @(3, 1, 72.47344466666667, 217.420334, 0.04973458223193756, 0.04973458223193756, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 34.57430230175445, 311.16872071579, 0.005045928987162869, 0.005045928987162868) },
    @{ Row = 14; A = "sCs"; B = "Gnas"; C = "Tshr"; D = "ECs"; Vals = @(3, 1, 343.5753833333333, 1030.72615, 0.2357770937183061, 0.235777093718306, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 247.7126977525834, 2229.41427977325, 0.03615230384604601, 0.036152303846046) },
    @{ Row = 15; A = "sCs"; B = "Gnas"; C = "Tshr"; D = "FAPs"; Vals = @(3, 1, 343.5753833333333, 1030.72615, 0.2357770937183061, 0.235777093718306, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 1163.794155741405, 10474.14740167265, 0.1698493469020286, 0.1698493469020285) },
    @{ Row = 16; A = "sCs"; B = "Gnas"; C = "Tshr"; D = "M2"; Vals = @(3, 1, 343.5753833333333, 1030.72615, 0.2357770937183061, 0.235777093718306, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 40.11231147903889, 361.01080331135, 0.005854170923469868, 0.005854170923469866) },
    @{ Row = 17; A = "sCs"; B = "Gnas"; C = "Tshr"; D = "sCs"; Vals = @(3, 1, 343.5753833333333, 1030.72615, 0.2357770937183061, 0.235777093718306, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 163.9066449986389, 1475.15980498775, 0.02392127204676166, 0.02392127204676165) }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    for ($i = 0; $i -lt $numCols.Length; $i++) {
        $ws.Range("$($numCols[$i])$r").Value = $item.Vals[$i]
    }
}
